# Update FuelPrices at 2025-03-31 07:55
# - Columns A (MLBSO00) and C (LNBSF00) are swapped (header + all data rows)
# - The "last row" date formatting (currently on B14) moves to the new last row B15
# - A new data row (15) is appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 14

# --- 1. Swap columns A and C for the header row and all existing data rows ---
for ($r = 1; $r -le $lastDataRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellC = $ws.Cells.Item($r, 3)
    $valA = $cellA.Value2
    $valC = $cellC.Value2
    $cellA.Value = $valC
    $cellC.Value = $valA
}

# --- 2. Remember the special "last row" date format currently on B14, then
#        normalize B14 back to the regular date format used by the other rows ---
$lastRowDateFormat = $ws.Range("B14").NumberFormat
$ws.Range("B14").NumberFormat = $ws.Range("B13").NumberFormat

# --- 3. Append the new data row (15) with the swapped layout already in effect ---
$newRow = 15
$ws.Cells.Item($newRow, 1).Value = 760.188
$ws.Cells.Item($newRow, 2).Value = 45744
$ws.Cells.Item($newRow, 3).Value = 800.9299999999999

# Give the new last row the special date format that used to belong to row 14
$ws.Range("B15").NumberFormat = $lastRowDateFormat
